$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.Formula = "'66.049.24"
$c.Style = "Normal"
$c = $ws.Cells.Item(2, 5)
$c.Formula = "'  -5.11%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(3, 4)
$c.Formula = "'3.286.84"
$c.Style = "Normal"
$c = $ws.Cells.Item(3, 5)
$c.Formula = "'  -5.59%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(4, 5)
$c.Formula = "'  -0.11%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(5, 4)
$c.Formula = "'559.65"
$c.Style = "Normal"
$c = $ws.Cells.Item(5, 5)
$c.Formula = "'  -3.42%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(6, 4)
$c.Formula = "'185.52"
$c.Style = "Normal"
$c = $ws.Cells.Item(6, 5)
$c.Formula = "'  -3.70%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(7, 5)
$c.Formula = "'  +0.08%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(8, 4)
$c.Formula = "'0.593"
$c.Style = "Normal"
$c = $ws.Cells.Item(8, 5)
$c.Formula = "'  -3.00%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(9, 4)
$c.Formula = "'3.280.32"
$c.Style = "Normal"
$c = $ws.Cells.Item(9, 5)
$c.Formula = "'  -5.45%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(10, 5)
$c.Formula = "'  -8.80%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(11, 4)
$c.Formula = "'0.586"
$c.Style = "Normal"
$c = $ws.Cells.Item(11, 5)
$c.Formula = "'  -5.16%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(12, 5)
$c.Formula = "'  -7.66%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(13, 4)
$c.Formula = "'0.0000267"
$c.Style = "Normal"
$c = $ws.Cells.Item(13, 5)
$c.Formula = "'  -6.83%  "
$c.Style = "Normal"
$ws.Cells.Item(14, 2).Value2 = "Polkadot"
$ws.Cells.Item(14, 3).Value2 = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$c = $ws.Cells.Item(14, 4)
$c.Formula = "'8.64"
$c.Style = "Normal"
$c = $ws.Cells.Item(14, 5)
$c.Formula = "'  -5.46%  "
$c.Style = "Normal"
$ws.Cells.Item(15, 2).Value2 = "BitcoinCash"
$ws.Cells.Item(15, 3).Value2 = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$c = $ws.Cells.Item(15, 4)
$c.Formula = "'632.73"
$c.Style = "Normal"
$c = $ws.Cells.Item(15, 5)
$c.Formula = "'  -2.44%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(16, 4)
$c.Formula = "'3.805.77"
$c.Style = "Normal"
$c = $ws.Cells.Item(16, 5)
$c.Formula = "'  -5.78%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(17, 4)
$c.Formula = "'66.046.38"
$c.Style = "Normal"
$c = $ws.Cells.Item(17, 5)
$c.Formula = "'  -4.86%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(18, 4)
$c.Formula = "'18.03"
$c.Style = "Normal"
$c = $ws.Cells.Item(18, 5)
$c.Formula = "'  -0.96%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(19, 4)
$c.Formula = "'0.117"
$c.Style = "Normal"
$c = $ws.Cells.Item(19, 5)
$c.Formula = "'  -3.32%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(20, 4)
$c.Formula = "'3.278.83"
$c.Style = "Normal"
$c = $ws.Cells.Item(20, 5)
$c.Formula = "'  -5.70%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(21, 4)
$c.Formula = "'11.37"
$c.Style = "Normal"
$c = $ws.Cells.Item(21, 5)
$c.Formula = "'  -7.87%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(22, 4)
$c.Formula = "'0.908"
$c.Style = "Normal"
$c = $ws.Cells.Item(22, 5)
$c.Formula = "'  -4.13%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(23, 4)
$c.Formula = "'18.20"
$c.Style = "Normal"
$c = $ws.Cells.Item(23, 5)
$c.Formula = "'  +1.81%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(24, 4)
$c.Formula = "'107.00"
$c.Style = "Normal"
$c = $ws.Cells.Item(24, 5)
$c.Formula = "'  +7.85%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(25, 4)
$c.Formula = "'4.90"
$c.Style = "Normal"
$c = $ws.Cells.Item(25, 5)
$c.Formula = "'  -7.05%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(26, 4)
$c.Formula = "'3.96"
$c.Style = "Normal"
$c = $ws.Cells.Item(26, 5)
$c.Formula = "'  -7.63%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(27, 5)
$c.Formula = "'  -7.08%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(28, 4)
$c.Formula = "'9.58"
$c.Style = "Normal"
$c = $ws.Cells.Item(28, 5)
$c.Formula = "'  -4.32%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(29, 4)
$c.Formula = "'8.72"
$c.Style = "Normal"
$c = $ws.Cells.Item(29, 5)
$c.Formula = "'  -6.56%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(30, 4)
$c.Formula = "'30.38"
$c.Style = "Normal"
$c = $ws.Cells.Item(30, 5)
$c.Formula = "'  -6.61%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(31, 4)
$c.Formula = "'4.03"
$c.Style = "Normal"
$c = $ws.Cells.Item(31, 5)
$c.Formula = "'  -6.38%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(32, 4)
$c.Formula = "'6.28"
$c.Style = "Normal"
$c = $ws.Cells.Item(32, 5)
$c.Formula = "'  -6.63%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(33, 4)
$c.Formula = "'11.07"
$c.Style = "Normal"
$c = $ws.Cells.Item(33, 5)
$c.Formula = "'  -4.88%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(34, 5)
$c.Formula = "'  -3.88%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(35, 4)
$c.Formula = "'57.70"
$c.Style = "Normal"
$c = $ws.Cells.Item(35, 5)
$c.Formula = "'  -5.31%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(36, 4)
$c.Formula = "'531.03"
$c.Style = "Normal"
$c = $ws.Cells.Item(36, 5)
$c.Formula = "'  +0.61%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(37, 4)
$c.Formula = "'3.719.48"
$c.Style = "Normal"
$c = $ws.Cells.Item(37, 5)
$c.Formula = "'  -0.81%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(38, 4)
$c.Formula = "'0.999"
$c.Style = "Normal"
$c = $ws.Cells.Item(38, 5)
$c.Formula = "'  -0.15%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(39, 4)
$c.Formula = "'3.35"
$c.Style = "Normal"
$c = $ws.Cells.Item(39, 5)
$c.Formula = "'  -4.81%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(40, 5)
$c.Formula = "'  -8.39%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(41, 4)
$c.Formula = "'2.74"
$c.Style = "Normal"
$c = $ws.Cells.Item(41, 5)
$c.Formula = "'  -6.91%  "
$c.Style = "Normal"
$ws.Cells.Item(42, 2).Value2 = "Kaspa"
$ws.Cells.Item(42, 3).Value2 = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Cells.Item(42, 4)
$c.Formula = "'0.129"
$c.Style = "Normal"
$c = $ws.Cells.Item(42, 5)
$c.Formula = "'  -3.34%  "
$c.Style = "Normal"
$ws.Cells.Item(43, 2).Value2 = "CoreDAO"
$ws.Cells.Item(43, 3).Value2 = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$c = $ws.Cells.Item(43, 4)
$c.Formula = "'3.41"
$c.Style = "Normal"
$c = $ws.Cells.Item(43, 5)
$c.Formula = "'  -2.55%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(44, 4)
$c.Formula = "'33.03"
$c.Style = "Normal"
$c = $ws.Cells.Item(44, 5)
$c.Formula = "'  -3.91%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(45, 4)
$c.Formula = "'0.339"
$c.Style = "Normal"
$c = $ws.Cells.Item(45, 5)
$c.Formula = "'  -9.21%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(46, 4)
$c.Formula = "'3.27"
$c.Style = "Normal"
$c = $ws.Cells.Item(46, 5)
$c.Formula = "'  -2.30%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(47, 5)
$c.Formula = "'  -5.95%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(48, 5)
$c.Formula = "'  -3.80%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(49, 5)
$c.Formula = "'  -7.81%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(50, 5)
$c.Formula = "'  -0.10%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(51, 4)
$c.Formula = "'1.27"
$c.Style = "Normal"
$c = $ws.Cells.Item(51, 5)
$c.Formula = "'  +3.02%  "
$c.Style = "Normal"
